$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.884.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.52%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.144.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.30%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "201.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "624.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.65%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.219"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.51%  "

$ws.Range("E9").Value = "  +1.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.488"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.81%  "

$ws.Range("E11").Value = "  +0.61%  "

$ws.Range("E12").Value = "  +6.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.724.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.22%  "

$ws.Range("E15").Value = "  +7.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.775.42"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.140.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +21.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "404.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.63"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.55"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.67"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.301.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "74.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("E29").Value = "  +4.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.994"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "521.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.54%  "

$ws.Range("E34").Value = "  +7.47%  "

$ws.Range("E35").Value = "  +21.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "21.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.59%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "196.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.38%  "

$ws.Range("E40").Value = "  -1.59%  "

$ws.Range("E41").Value = "  +0.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.104"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.22%  "

$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.813"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +21.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.65%  "

$ws.Range("E49").Value = "  +7.49%  "

$ws.Range("E50").Value = "  +4.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.86%  "
